$d = $word.ActiveDocument

# --- Change 1: update the timestamp in the document date line ---
[void]$d.Content.Find.Execute(
    "June  30, 2021 (04:23:45 PM)", $true, $false, $false, $false, $false,
    $true, 1, $false, "June  30, 2021 (05:30:20 PM)", 2)

# --- Change 2: expand "it should have two attributes, of type" into
#     "it should have two attributes, width and length, of type" with
#     "width" and "length" styled as NormalTok (code-token) runs. ---

# Locate the specific list-item paragraph so the Find below only ever
# touches this occurrence (the words "width"/"length" also appear,
# already styled, elsewhere in the document).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*it should have two attributes, of type*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $rng = $target.Range
    [void]$rng.Find.Execute(
        "it should have two attributes, of type", $true, $false, $false,
        $false, $false, $true, 1, $false,
        "it should have two attributes, width and length, of type", 2)

    # Style "width" within that same paragraph.
    $rngWidth = $target.Range
    [void]$rngWidth.Find.Execute("width", $true, $false, $false, $false,
        $false, $true, 1, $false, "", 0)
    $rngWidth.Style = "NormalTok"

    # Style "length" within that same paragraph.
    $rngLength = $target.Range
    [void]$rngLength.Find.Execute("length", $true, $false, $false, $false,
        $false, $true, 1, $false, "", 0)
    $rngLength.Style = "NormalTok"
}
